$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = [double]"9.532633749233283E-05"
$ws.Range("E2").Value = [double]"9.532633749233283E-05"

$ws.Range("D3").Value = [double]"0.8858080328161515"
$ws.Range("E3").Value = [double]"0.8858080328161515"

$ws.Range("D4").Value = [double]"7.654298864016087E-05"
$ws.Range("E4").Value = [double]"7.654298864016087E-05"

$ws.Range("D5").Value = [double]"1.645060661414587E-15"
$ws.Range("E5").Value = [double]"1.645060661414587E-15"

$ws.Range("D6").Value = [double]"0.5789041841298012"
$ws.Range("E6").Value = [double]"0.5789041841298012"

$ws.Range("C7").Value = $false
$ws.Range("D7").Value = [double]"1.605788219369724E-07"
$ws.Range("E7").Value = [double]"0.9999998394211781"

$ws.Range("D8").Value = [double]"0.9302801168955777"
$ws.Range("E8").Value = [double]"0.06971988310442234"

$ws.Range("D9").Value = [double]"0.9823815057254937"
$ws.Range("E9").Value = [double]"0.01761849427450635"

$ws.Range("D10").Value = [double]"0.9999999164888559"
$ws.Range("E10").Value = [double]"8.351114411109251E-08"

$ws.Range("D11").Value = [double]"1"
$ws.Range("E11").Value = [double]"0"
$ws.Range("F11").Value = [double]"1.876946687698364"
$ws.Range("G11").Value = [double]"0.7"

$ws.Range("D12").Value = [double]"2.601121541647324E-05"
$ws.Range("E12").Value = [double]"2.601121541647324E-05"

$ws.Range("D13").Value = [double]"0.9734773640976235"
$ws.Range("E13").Value = [double]"0.9734773640976235"

$ws.Range("D14").Value = [double]"7.267181660675463E-06"
$ws.Range("E14").Value = [double]"7.267181660675463E-06"

$ws.Range("D15").Value = [double]"2.821612414598692E-26"
$ws.Range("E15").Value = [double]"2.821612414598692E-26"

$ws.Range("D16").Value = [double]"0.3698691962212504"
$ws.Range("E16").Value = [double]"0.3698691962212504"

$ws.Range("C17").Value = $false
$ws.Range("D17").Value = [double]"1.805948999449052E-12"
$ws.Range("E17").Value = [double]"0.999999999998194"

$ws.Range("D18").Value = [double]"0.9856250337970148"
$ws.Range("E18").Value = [double]"0.01437496620298517"

$ws.Range("D19").Value = [double]"0.9653291138225722"
$ws.Range("E19").Value = [double]"0.03467088617742775"

$ws.Range("D20").Value = [double]"0.9999999999998679"
$ws.Range("E20").Value = [double]"1.321165399303936E-13"

$ws.Range("D21").Value = [double]"1"
$ws.Range("E21").Value = [double]"0"
$ws.Range("F21").Value = [double]"3.118131875991821"
$ws.Range("G21").Value = [double]"0.8"
